$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.783.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.38%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.92"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.27%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.441.83"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.23%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.49%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.79%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.57%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.89"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.877.26"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.283.92"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.439.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.84"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.81%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +14.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.76"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "617.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.88%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.37"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.52%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.541.10"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.52%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.993"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.67%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.46%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.51%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.39%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.29%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.73%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.90%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.02%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +21.65%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +7.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.32"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.07%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0281"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.48"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.68%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.59"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.26"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.61%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.41%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.49%  "
